$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor style (default/unstyled) used to restore D-column cell style after
# forcing a Text number-format, so numeric-looking price strings (e.g. "1.000")
# are stored as text instead of being auto-converted to numbers by Excel.
$donorStyle = $ws.Range("D6").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.413.93"
$ws.Range("D2").Style = $donorStyle
$ws.Range("E2").Value = "  +1.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.823.50"
$ws.Range("D3").Style = $donorStyle
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = $donorStyle
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.05"
$ws.Range("D5").Style = $donorStyle
$ws.Range("E5").Value = "  +0.70%  "

$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4487"
$ws.Range("D7").Style = $donorStyle
$ws.Range("E7").Value = "  +1.85%  "

$ws.Range("E8").Value = "  +2.10%  "

$ws.Range("E9").Value = "  +3.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8864"
$ws.Range("D10").Style = $donorStyle
$ws.Range("E10").Value = "  +4.96%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.16"
$ws.Range("D11").Style = $donorStyle
$ws.Range("E11").Value = "  +2.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.827.13"
$ws.Range("D12").Style = $donorStyle
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.763"
$ws.Range("D13").Style = $donorStyle
$ws.Range("E13").Value = "  +1.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.06"
$ws.Range("D14").Style = $donorStyle
$ws.Range("E14").Value = "  +4.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.412"
$ws.Range("D15").Style = $donorStyle
$ws.Range("E15").Value = "  +2.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07097"
$ws.Range("D16").Style = $donorStyle
$ws.Range("E16").Value = "  +0.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("D17").Style = $donorStyle
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008800"
$ws.Range("D18").Style = $donorStyle
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.0000"
$ws.Range("D19").Style = $donorStyle
$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.18"
$ws.Range("D20").Style = $donorStyle
$ws.Range("E20").Value = "  +1.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.411.11"
$ws.Range("D21").Style = $donorStyle
$ws.Range("E21").Value = "  +1.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.299"
$ws.Range("D22").Style = $donorStyle
$ws.Range("E22").Value = "  +2.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.94"
$ws.Range("D23").Style = $donorStyle
$ws.Range("E23").Value = "  +0.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.056.58"
$ws.Range("D24").Style = $donorStyle
$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.958"
$ws.Range("D25").Style = $donorStyle
$ws.Range("E25").Value = "  -1.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.366"
$ws.Range("D26").Style = $donorStyle
$ws.Range("E26").Value = "  +7.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.54"
$ws.Range("D27").Style = $donorStyle
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  +1.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.373"
$ws.Range("D29").Style = $donorStyle
$ws.Range("E29").Value = "  +2.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.10"
$ws.Range("D30").Style = $donorStyle
$ws.Range("E30").Value = "  +0.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08851"
$ws.Range("D31").Style = $donorStyle
$ws.Range("E31").Value = "  +0.80%  "

$ws.Range("E32").Value = "  +6.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.202"
$ws.Range("D33").Style = $donorStyle
$ws.Range("E33").Value = "  +2.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.530"
$ws.Range("D34").Style = $donorStyle
$ws.Range("E34").Value = "  +2.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.921"
$ws.Range("D35").Style = $donorStyle
$ws.Range("E35").Value = "  +1.10%  "

$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.114"
$ws.Range("D37").Style = $donorStyle
$ws.Range("E37").Value = "  +1.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01998"
$ws.Range("D38").Style = $donorStyle
$ws.Range("E38").Value = "  +2.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05324"
$ws.Range("D39").Style = $donorStyle
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.371"
$ws.Range("D40").Style = $donorStyle
$ws.Range("E40").Value = "  +1.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5326"
$ws.Range("D41").Style = $donorStyle
$ws.Range("E41").Value = "  +3.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1723"
$ws.Range("D42").Style = $donorStyle
$ws.Range("E42").Value = "  +1.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.854"
$ws.Range("D43").Style = $donorStyle
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.305"
$ws.Range("D44").Style = $donorStyle
$ws.Range("E44").Value = "  +19.21%  "

$ws.Range("E45").Value = "  +2.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5102"
$ws.Range("D46").Style = $donorStyle
$ws.Range("E46").Value = "  +5.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.63"
$ws.Range("D47").Style = $donorStyle
$ws.Range("E47").Value = "  +0.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.703"
$ws.Range("D48").Style = $donorStyle
$ws.Range("E48").Value = "  +2.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.78"
$ws.Range("D49").Style = $donorStyle
$ws.Range("E49").Value = "  -0.24%  "

$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06386"
$ws.Range("D51").Style = $donorStyle
$ws.Range("E51").Value = "  +0.79%  "

